$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("B1").Value = "NO_PRED"
$ws.Range("C1").Value = "NEAR_PRED"
$ws.Range("D1").Value = "BEING_EATEN"

# NO_PRED block (IDLE row removed, remaining rows shifted up)
$ws.Range("A2").Value = "NO_PREDRAND"
$ws.Range("B2").Value = 0.99
$ws.Range("C2").Value = -0.5
$ws.Range("D2").Value = -0.99

$ws.Range("A3").Value = "NO_PREDHUNT"
$ws.Range("B3").Value = 0.99
$ws.Range("C3").Value = -0.5
$ws.Range("D3").Value = -0.99

$ws.Range("A4").Value = "NO_PREDMULTIPLY"
$ws.Range("B4").Value = 0.99
$ws.Range("C4").Value = -0.5
$ws.Range("D4").Value = -0.99

$ws.Range("A5").Value = "NO_PREDFLEE"
$ws.Range("B5").Value = -0.99
$ws.Range("C5").Value = -0.99
$ws.Range("D5").Value = -0.99

# NEAR_PRED block (IDLE row removed, remaining rows shifted up)
$ws.Range("A6").Value = "NEAR_PREDRAND"
$ws.Range("B6").Value = 0.99
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = -0.99

$ws.Range("A7").Value = "NEAR_PREDHUNT"
$ws.Range("B7").Value = 0.99
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.99

$ws.Range("A8").Value = "NEAR_PREDMULTIPLY"
$ws.Range("B8").Value = 0.99
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = -0.99

$ws.Range("A9").Value = "NEAR_PREDFLEE"
$ws.Range("B9").Value = 0.5
$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = -0.99

# BEING_EATEN block (IDLE row removed, remaining rows shifted up)
$ws.Range("A10").Value = "BEING_EATENRAND"
$ws.Range("B10").Value = 0.99
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = -0.99

$ws.Range("A11").Value = "BEING_EATENHUNT"
$ws.Range("B11").Value = 0.99
$ws.Range("C11").Value = 0.5
$ws.Range("D11").Value = -0.99

$ws.Range("A12").Value = "BEING_EATENMULTIPLY"
$ws.Range("B12").Value = 0.99
$ws.Range("C12").Value = 0.5
$ws.Range("D12").Value = -0.99

$ws.Range("A13").Value = "BEING_EATENFLEE"
$ws.Range("B13").Value = 0.99
$ws.Range("C13").Value = 0.99
$ws.Range("D13").Value = -0.99

# Remove the now-unused trailing rows (formerly BEING_EATENHUNT/MULTIPLY/FLEE)
$ws.Range("A14:D16").ClearContents()

# Restore the selection recorded in the saved workbook
$ws.Range("C17").Select() | Out-Null
